$wb = $excel.ActiveWorkbook

# --- Daily_Data: append a new date block (date serial 46029) of 22 rows ---
$ws = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @{ Name = "ASAHI DEPOSITORY LLC Registered"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 },
    @{ Name = "ASAHI DEPOSITORY LLC Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 },
    @{ Name = "BRINK'S, INC. Registered"; C = 90027.72500000001; D = 0; E = 0; F = 0; G = 0; H = 90027.72500000001 },
    @{ Name = "BRINK'S, INC. Eligible"; C = 5744.711; D = 0; E = 0; F = 0; G = 0; H = 5744.711 },
    @{ Name = "CNT DEPOSITORY, INC. Registered"; C = 1246.06; D = 0; E = 0; F = 0; G = 0; H = 1246.06 },
    @{ Name = "CNT DEPOSITORY, INC. Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 },
    @{ Name = "DELAWARE DEPOSITORY Registered"; C = 1633.941; D = 0; E = 0; F = 0; G = 0; H = 1633.941 },
    @{ Name = "DELAWARE DEPOSITORY Eligible"; C = 18509.729; D = 0; E = 0; F = 0; G = 0; H = 18509.729 },
    @{ Name = "HSBC BANK, USA Registered"; C = 1295.223; D = 0; E = 0; F = 0; G = 0; H = 1295.223 },
    @{ Name = "HSBC BANK, USA Eligible"; C = 9281.978999999999; D = 0; E = 0; F = 0; G = 0; H = 9281.978999999999 },
    @{ Name = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"; C = 2395.448; D = 0; E = 0; F = 0; G = 0; H = 2395.448 },
    @{ Name = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 },
    @{ Name = "JP MORGAN CHASE BANK NA Registered"; C = 124991.729; D = 0; E = 0; F = 0; G = 0; H = 124991.729 },
    @{ Name = "JP MORGAN CHASE BANK NA Eligible"; C = 125407.673; D = 0; E = 0; F = 0; G = 0; H = 125407.673 },
    @{ Name = "LOOMIS INTERNATIONAL (US) LLC Registered"; C = 68084.33; D = 0; E = 0; F = 0; G = 0; H = 68084.33 },
    @{ Name = "LOOMIS INTERNATIONAL (US) LLC Eligible"; C = 106188.481; D = 0; E = 0; F = 0; G = 0; H = 106188.481 },
    @{ Name = "MALCA-AMIT USA, LLC Registered"; C = 395.145; D = 0; E = 0; F = 0; G = 0; H = 395.145 },
    @{ Name = "MALCA-AMIT USA, LLC Eligible"; C = 0; D = 0; E = 0; F = 0; G = 0; H = 0 },
    @{ Name = "MANFRA, TORDELLA & BROOKES, LLC Registered"; C = 54605.27; D = 0; E = 0; F = 0; G = 0; H = 54605.27 },
    @{ Name = "MANFRA, TORDELLA & BROOKES, LLC Eligible"; C = 21419.744; D = 0; E = 20351.336; F = -20351.336; G = 0; H = 1068.408 },
    @{ Name = "STONEX PRECIOUS METALS LLC Registered"; C = 14122.765; D = 0; E = 0; F = 0; G = 0; H = 14122.765 },
    @{ Name = "STONEX PRECIOUS METALS LLC Eligible"; C = 16.075; D = 0; E = 0; F = 0; G = 0; H = 16.075 }
)

$startRow = 68
$dateSerial = 46029
$dateNumberFormat = $ws.Cells.Item($startRow - 1, 1).NumberFormat

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 1).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}

# --- Today_Summary: MANFRA, TORDELLA & BROOKES, LLC row (row 11) updated ---
$ws2 = $wb.Worksheets.Item("Today_Summary")
$ws2.Cells.Item(11, 2).Value = 1068.408
$ws2.Cells.Item(11, 4).Value = 55673.678

# --- Monthly_Stats: top summary row (row 2) and MANFRA Eligible detail row (row 25) ---
$ws3 = $wb.Worksheets.Item("Monthly_Stats")
$ws3.Cells.Item(2, 2).Value = 266217.056
$ws3.Cells.Item(2, 4).Value = 625014.692
$ws3.Cells.Item(25, 4).Value = 27826.734
$ws3.Cells.Item(25, 5).Value = 1068.408
